# Updated cryptos list on Sat Apr  6 05:31:35 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the coin rows on
# the active sheet, and swaps the TheGraph / ApeXProtocol rows (44 <-> 45)
# to reflect their new ranking order, with updated price/volume for both.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a "Price"-style cell without letting
# Excel's automatic type-inference turn a plain decimal-looking string
# (e.g. "582.83") into a numeric value. Values that contain thousands
# separators (two dots, e.g. "67.857.14"), subscript digits, or other
# non-numeric characters are safe to assign directly.
function Set-TextValue($cellRef, $value) {
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($cellRef).Value = "'" + $value
    } else {
        $ws.Range($cellRef).Value = $value
    }
}

# row -> (new Price, new Volume(1h)) ; $null means "unchanged"
$updates = @{
    2  = @("67.857.14", "  +2.08%  ")
    3  = @("3.339.05",  "  +2.73%  ")
    4  = @($null,       "  +0.07%  ")
    5  = @("582.83",    "  +3.17%  ")
    6  = @("176.82",    "  +1.74%  ")
    7  = @("1.00",      "  -0.22%  ")
    8  = @($null,       "  +2.10%  ")
    9  = @("3.333.82",  "  +2.70%  ")
    10 = @($null,       "  +5.78%  ")
    11 = @("0.580",     "  +2.85%  ")
    12 = @("46.85",     "  +4.38%  ")
    13 = @($null,       "  +2.79%  ")
    14 = @("690.44",    "  -0.39%  ")
    15 = @("3.887.21",  "  +2.99%  ")
    16 = @("8.45",      $null)
    17 = @("67.890.64", "  +1.83%  ")
    18 = @($null,       "  +0.00%  ")
    19 = @("3.337.04",  "  +2.17%  ")
    20 = @($null,       "  +1.57%  ")
    21 = @($null,       "  +4.69%  ")
    22 = @($null,       "  +2.06%  ")
    23 = @("5.37",      "  +6.10%  ")
    24 = @("17.03",     "  +1.74%  ")
    25 = @("98.60",     "  +1.59%  ")
    26 = @($null,       "  +1.57%  ")
    27 = @($null,       "  +0.81%  ")
    28 = @("9.53",      "  +3.92%  ")
    29 = @($null,       "  +1.68%  ")
    30 = @($null,       "  +3.28%  ")
    31 = @("7.10",      "  +7.58%  ")
    32 = @("572.96",    "  +0.00%  ")
    33 = @($null,       "  +3.24%  ")
    34 = @($null,       "  +3.70%  ")
    35 = @("3.721.30",  "  -2.73%  ")
    36 = @("57.30",     "  +3.98%  ")
    37 = @("0.999",     "  -0.02%  ")
    38 = @("3.34",      "  +3.40%  ")
    39 = @("34.24",     "  +9.07%  ")
    40 = @($null,       "  +2.99%  ")
    41 = @($null,       "  +5.07%  ")
    42 = @($null,       "  +3.28%  ")
    43 = @("0.0₃0676",  "  +2.55%  ")
    46 = @("0.0406",    "  +1.41%  ")
    47 = @($null,       "  +7.31%  ")
    48 = @($null,       "  +2.15%  ")
    49 = @($null,       "  -0.58%  ")
    50 = @($null,       "  -2.32%  ")
    51 = @("129.51",    "  +0.86%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]
    if ($price -ne $null) {
        Set-TextValue "D$row" $price
    }
    if ($volume -ne $null) {
        $ws.Range("E$row").Value = $volume
    }
}

# TheGraph and ApeXProtocol swap ranking places (row 44 <-> row 45), each
# carrying its own refreshed price/volume figures.
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D44" "3.33"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D45" "0.336"
$ws.Range("E45").Value = "  +4.03%  "
